$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.031.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.907.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4639'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.37%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4065'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.72%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.90'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08000'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.005'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.916.12'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.940'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("E15").Value = '  -1.95%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.47%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001035'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06542'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.012.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.467'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.26%  '

$ws.Range("E24").Value = '  +1.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.246'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.132.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.105'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.404'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9802'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09381'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.419'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("E35").Value = '  +0.79%  '

$ws.Range("E36").Value = '  -0.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06081'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02229'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.401'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.162'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5814'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9996'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.73%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1824'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.258'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.359'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +15.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5488'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.903'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07027'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +22.31%  '
